$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A8").Value = "LD48"
$ws.Range("B8").Value = "Bouldaouf dash"
$ws.Range("C8").Value = 3.235
$ws.Range("D8").Value = 3.103
$ws.Range("E8").Value = 2.765
$ws.Range("K1").Formula = "=""2.603"""
$ws.Range("K1").Copy()
$ws.Range("F8").PasteSpecial(-4163)
$ws.Range("K1").Formula = "=""2.735"""
$ws.Range("K1").Copy()
$ws.Range("G8").PasteSpecial(-4163)
$ws.Range("K1").ClearContents()
$excel.CutCopyMode = $false
$ws.Range("H8").Value = 2.838
$ws.Range("I8").Value = 2.758

$ws.Range("B8").Select()
